$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Petalburg Woods grunt (TRAINER_GRUNT_PETALBURG_WOODS): fix sprite/class ---
# Insert 3 rows right after the existing "Charmander" party line (row 156) and
# before "TRAINER_LYLE" (row 158) to add trainerClass / encounterMusic_gender /
# trainerPic overrides so the grunt uses the Team Magma grunt (female) sprite.
$ws.Rows("157:159").Insert()

$ws.Range("A157").Value = ".trainerClass"
$ws.Range("B157").Value = "TRAINER_CLASS_TEAM_MAGMA"

$ws.Range("A158").Value = ".encounterMusic_gender"
$ws.Range("B158").Value = "TRAINER_ENCOUNTER_MUSIC_MAGMA"

$ws.Range("A159").Value = ".trainerPic"
$ws.Range("B159").Value = "TRAINER_PIC_MAGMA_GRUNT_F"

# --- New trainer block before END (second TRAINER_HALEY_1 entry) ---
# After the previous insertion, the old "END" row (190) now sits at row 193.
# Insert 5 rows there to make room for a new trainer section, pushing END
# further down to row 198.
$ws.Rows("193:197").Insert()

$ws.Range("A193").Value = "TRAINER_HALEY_1"

$ws.Range("A194").Value = "species"
$ws.Range("B194").Value = "lvl"
$ws.Range("C194").Value = "iv"
$ws.Range("D194").Value = "heldItem"
$ws.Range("E194").Value = "moves"

$ws.Range("A195").Value = "Geodude"
$ws.Range("B195").Value = 5

$ws.Range("A196").Value = "Anorith"
$ws.Range("B196").Value = 6

# --- Scroll/selection: move view toward Route 103 area ---
try {
    $win = $excel.ActiveWindow
    $win.ScrollRow = 74
    $win.ScrollColumn = 1
} catch {
}
$null = $ws.Range("A74").Select()
